$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 504
$ws.Range("J2").Value = 724.6667
$ws.Range("L2").Value = 724.6667
$ws.Range("N2").Value = -950.6667

$ws.Range("H4").Value = 164.25
$ws.Range("I4").Value = 164.25
$ws.Range("K4").Value = 164.25
$ws.Range("M4").Value = -50.25

$ws.Range("H9").Value = 93.17391000000001
$ws.Range("I9").Value = 98.52381
$ws.Range("J9").Value = 37
$ws.Range("K9").Value = 98.52381
$ws.Range("L9").Value = 37
$ws.Range("M9").Value = 70.47619
$ws.Range("N9").Value = -375

$ws.Range("H19").Value = 2476
$ws.Range("I19").Value = 2286.077
$ws.Range("J19").Value = 2700.4546
$ws.Range("K19").Value = 2286.077
$ws.Range("L19").Value = 2700.4546
$ws.Range("M19").Value = -2111.077
$ws.Range("N19").Value = -3050.4546

$ws.Range("H70").Value = 73062.36
$ws.Range("J70").Value = 112603.664
$ws.Range("L70").Value = 337810.992
$ws.Range("N70").Value = -338350.992

$ws.Range("H73").Value = 73062.36
$ws.Range("J73").Value = 112603.664
$ws.Range("L73").Value = 337810.992
$ws.Range("N73").Value = -339682.992

$ws.Range("H101").Value = 621.3333
$ws.Range("I101").Value = 548.4
$ws.Range("K101").Value = 1645.2
$ws.Range("M101").Value = -23.19999999999982

$ws.Range("H112").Value = 3231.524
$ws.Range("J112").Value = 3534.3333
$ws.Range("L112").Value = 10602.9999
$ws.Range("N112").Value = -12818.9999

$ws.Range("H116").Value = 4986.25
$ws.Range("I116").Value = 4978.8
$ws.Range("K116").Value = 4978.8
$ws.Range("M116").Value = -1536.8

$ws.Range("H137").Value = 3646.8333
$ws.Range("I137").Value = 3577.3333
$ws.Range("K137").Value = 10731.9999
$ws.Range("M137").Value = -8181.999899999999

$ws.Range("H138").Value = 3256.5278
$ws.Range("I138").Value = 1696.3636
$ws.Range("J138").Value = 3943
$ws.Range("K138").Value = 5089.0908
$ws.Range("L138").Value = 11829
$ws.Range("M138").Value = 50.90920000000006
$ws.Range("N138").Value = -22109

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 8000
$ws.Range("J27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("N27").Value = -8368

$ws.Range("H32").Value = 10595.352
$ws.Range("I32").Value = 786.9286
$ws.Range("K32").Value = 786.9286
$ws.Range("M32").Value = -499.9286

$ws.Range("H74").Value = 3118.8462
$ws.Range("I74").Value = 2269.111
$ws.Range("J74").Value = 5030.75
$ws.Range("K74").Value = 2269.111
$ws.Range("L74").Value = 5030.75
$ws.Range("M74").Value = -1395.111
$ws.Range("N74").Value = -6778.75

$ws.Range("H77").Value = 3118.8462
$ws.Range("I77").Value = 2269.111
$ws.Range("J77").Value = 5030.75
$ws.Range("K77").Value = 11345.555
$ws.Range("L77").Value = 25153.75
$ws.Range("M77").Value = -6977.555
$ws.Range("N77").Value = -33889.75

$ws.Range("H95").Value = 23231
$ws.Range("J95").Value = 23231
$ws.Range("L95").Value = 23231
$ws.Range("N95").Value = -28723

$ws.Range("H102").Value = 1296.5834
$ws.Range("I102").Value = 1220.1052
$ws.Range("K102").Value = 1220.1052
$ws.Range("M102").Value = 401.8948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1888.9166
$ws.Range("I64").Value = 1914
$ws.Range("K64").Value = 1914
$ws.Range("M64").Value = -1689

$ws.Range("H67").Value = 1888.9166
$ws.Range("I67").Value = 1914
$ws.Range("K67").Value = 1914
$ws.Range("M67").Value = -1134

$ws.Range("H100").Value = 44412
$ws.Range("J100").Value = 44412
$ws.Range("L100").Value = 44412
$ws.Range("N100").Value = -46576

$ws.Range("H107").Value = 873.7143
$ws.Range("I107").Value = 873.7143
$ws.Range("K107").Value = 873.7143
$ws.Range("M107").Value = 1046.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 60227.695
$ws.Range("I16").Value = 32639.715
$ws.Range("K16").Value = 32639.715
$ws.Range("M16").Value = -32352.715

$ws.Range("H86").Value = 8070.5713
$ws.Range("I86").Value = 7498.8
$ws.Range("K86").Value = 7498.8
$ws.Range("M86").Value = -6375.8

$ws.Range("H89").Value = 8070.5713
$ws.Range("I89").Value = 7498.8
$ws.Range("K89").Value = 37494
$ws.Range("M89").Value = -31878

$ws.Range("H96").Value = 30936.889
$ws.Range("J96").Value = 30936.889
$ws.Range("L96").Value = 30936.889
$ws.Range("N96").Value = -36428.889

$ws.Range("H105").Value = 2880.3333
$ws.Range("I105").Value = 2760.6667
$ws.Range("K105").Value = 2760.6667
$ws.Range("M105").Value = -1013.6667

$ws.Range("H113").Value = 60227.695
$ws.Range("I113").Value = 32639.715
$ws.Range("K113").Value = 32639.715
$ws.Range("M113").Value = -30469.715

$ws.Range("H137").Value = 43077.7
$ws.Range("J137").Value = 40666.332
$ws.Range("L137").Value = 40666.332
$ws.Range("N137").Value = -50866.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 47
$ws.Range("J7").Value = 93.5
$ws.Range("K7").Value = 141
$ws.Range("L7").Value = 280.5
$ws.Range("M7").Value = -29
$ws.Range("N7").Value = -504.5

$ws.Range("H12").Value = 322.0435
$ws.Range("I12").Value = 373.83334
$ws.Range("K12").Value = 1121.50002
$ws.Range("M12").Value = -948.5000199999999

$ws.Range("H33").Value = 91
$ws.Range("I33").Value = 91
$ws.Range("K33").Value = 546
$ws.Range("M33").Value = -263

$ws.Range("H129").Value = 2855.8333
$ws.Range("I129").Value = 601.6667
$ws.Range("J129").Value = 5110
$ws.Range("K129").Value = 1805.0001
$ws.Range("L129").Value = 15330
$ws.Range("M129").Value = 3194.9999
$ws.Range("N129").Value = -25330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 8999
$ws.Range("J23").Value = 8999
$ws.Range("L23").Value = 8999
$ws.Range("N23").Value = -9445

$ws.Range("H28").Value = 2000
$ws.Range("J28").Value = 2000
$ws.Range("L28").Value = 2000
$ws.Range("N28").Value = -2384

$ws.Range("H70").Value = 6002.6
$ws.Range("I70").Value = 5965
$ws.Range("J70").Value = 6106
$ws.Range("K70").Value = 5965
$ws.Range("L70").Value = 6106
$ws.Range("M70").Value = -5695
$ws.Range("N70").Value = -6646

$ws.Range("H73").Value = 6002.6
$ws.Range("I73").Value = 5965
$ws.Range("J73").Value = 6106
$ws.Range("K73").Value = 5965
$ws.Range("L73").Value = 6106
$ws.Range("M73").Value = -5029
$ws.Range("N73").Value = -7978

$ws.Range("H102").Value = 2971.9092
$ws.Range("J102").Value = 5950
$ws.Range("L102").Value = 5950
$ws.Range("N102").Value = -9194

$ws.Range("H106").Value = 45928.668
$ws.Range("J106").Value = 45928.668
$ws.Range("L106").Value = 45928.668
$ws.Range("N106").Value = -48452.668

$ws.Range("H122").Value = 4406.9473
$ws.Range("I122").Value = 2461.7144
$ws.Range("K122").Value = 7385.1432
$ws.Range("M122").Value = -4935.1432

$ws.Range("H132").Value = 3063.1052
$ws.Range("I132").Value = 3122.8333
$ws.Range("J132").Value = 1988
$ws.Range("K132").Value = 9368.499899999999
$ws.Range("L132").Value = 5964
$ws.Range("M132").Value = -6838.499899999999
$ws.Range("N132").Value = -11024

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6805.143
$ws.Range("I9").Value = 4109.5
$ws.Range("J9").Value = 10399.333
$ws.Range("K9").Value = 4109.5
$ws.Range("L9").Value = 10399.333
$ws.Range("M9").Value = -3885.5
$ws.Range("N9").Value = -10847.333

$ws.Range("H16").Value = 2029
$ws.Range("I16").Value = 1104.909
$ws.Range("J16").Value = 5417.3335
$ws.Range("K16").Value = 1104.909
$ws.Range("L16").Value = 5417.3335
$ws.Range("M16").Value = -934.9090000000001
$ws.Range("N16").Value = -5757.3335

$ws.Range("H82").Value = 1185.2354
$ws.Range("I82").Value = 1040.5625
$ws.Range("K82").Value = 1040.5625
$ws.Range("M82").Value = -679.5625

$ws.Range("H85").Value = 1185.2354
$ws.Range("I85").Value = 1040.5625
$ws.Range("K85").Value = 1040.5625
$ws.Range("M85").Value = 207.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1293.159
$ws.Range("I14").Value = 1090.6744
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 1090.6744
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -922.6744000000001
$ws.Range("N14").Value = -10336

$ws.Range("H113").Value = 651.8
$ws.Range("I113").Value = 660.9231
$ws.Range("J113").Value = 592.5
$ws.Range("K113").Value = 1982.7693
$ws.Range("L113").Value = 1777.5
$ws.Range("M113").Value = 187.2307000000001
